# Inserimento file input RFCF_fatturato in snapshot
# Adds 16 new rows (56-71) to the "Snapshot" sheet describing the new
# rfcf_fatturato / data / input tables used to feed the snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Copy the cell formatting (borders/fill/font) of the existing
#        data rows (A2:H2 uses style "3", no extra columns) down onto the
#        16 new rows A56:H71 so the new rows look like the rest of the
#        table. PasteSpecial(-4122) = xlPasteFormats, values are left
#        untouched (we set them explicitly afterwards).
$ws.Range("A2:H2").Copy() | Out-Null
for ($r = 56; $r -le 71; $r++) {
    $ws.Range("A$r`:H$r").PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = 0

# --- 2. Seed the "template" row first (row 62, rfcf_configurazione_start_run)
#        so that the new shared-string entries are created in the same
#        order the workbook was originally authored.
$ws.Range("A62").Value = "rfcf_fatturato"
$ws.Range("B62").Value = "data"
$ws.Range("C62").Value = "rfcf_configurazione_start_run"
$ws.Range("D62").Value = "input"

# --- 3. Fill in columns A (flusso), B (database) and D (tipo_storicizzazione)
#        for the remaining new rows - these reuse the strings just created
#        above, so no new shared-string entries are added here.
$otherRows = 56,57,58,59,60,61,63,64,65,66,67,68,69,70,71
foreach ($r in $otherRows) {
    $ws.Range("A$r").Value = "rfcf_fatturato"
    $ws.Range("B$r").Value = "data"
    $ws.Range("D$r").Value = "input"
}

# --- 4. Fill in column C (nome_tabella) for each new row, in row order,
#        skipping row 62 (already set above).
$ws.Range("C56").Value = "tcr_configurazione_storico_snapshot"
$ws.Range("C57").Value = "tcr_decodifiche"
$ws.Range("C58").Value = "tcr_codifiche_tariffario"
$ws.Range("C59").Value = "tcr_mapping_xe_gas"
$ws.Range("C60").Value = "tcr_mapping_xe_pwr"
$ws.Range("C61").Value = "tcr_tariffario_xe_indici_fatturazione"
$ws.Range("C63").Value = "rfcf_calendario_fatturazione_mensile_all_ver"
$ws.Range("C64").Value = "rfcf_calendario_xe_date"
$ws.Range("C65").Value = "rfcf_calendario_xe_forn"
$ws.Range("C66").Value = "rfcf_estrazione_CE"
$ws.Range("C67").Value = "rfcf_estrazione_DR"
$ws.Range("C68").Value = "rfcf_parametri_previsione_extra"
$ws.Range("C69").Value = "rfcf_forzatura_iva"
$ws.Range("C70").Value = "rfcf_ordine_default_tioce"
$ws.Range("C71").Value = "rfcf_tipo_regime"

# --- 5. Restore the view state (zoom + selection) to match the saved
#        workbook - scrolled down to the new rows, zoomed to 70%.
$win = $excel.ActiveWindow
$win.Zoom = 70
$ws.Range("C32").Select() | Out-Null

Write-Host "Added rows 56-71 to Snapshot sheet."
